# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values for the
# "8b95212a-cc2f-4185-9ddb-738e68c91732..." row (row 3) on both the
# zh-cn and de-de language sheets, reflecting a fresh report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 22:51:40"
$wsZhCn.Range("H3").Value = "2016-03-22 22:52:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 22:51:45"
$wsDeDe.Range("H3").Value = "2016-03-22 22:52:16"
